# Update the 25 division problems/answers in the single table on the page.
# Each cell is addressed directly by (row, column) and its text is replaced
# in place (rather than using Find/Replace across the whole document) so that
# duplicate values elsewhere in the table can never cause a wrong cell to be
# updated. The trailing cell-mark character of each cell's Range is excluded
# before assigning the new text so the table structure/formatting is kept.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1,1)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "17÷8=2, 1"

$cell = $t.Cell(1,2)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "30÷5=6, 0"

$cell = $t.Cell(1,3)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "14÷7=2, 0"

$cell = $t.Cell(1,4)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "31÷4=7, 3"

$cell = $t.Cell(1,5)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "89÷7=12, 5"

$cell = $t.Cell(5,1)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "59÷9=6, 5"

$cell = $t.Cell(5,2)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "48÷5=9, 3"

$cell = $t.Cell(5,3)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "22÷9=2, 4"

$cell = $t.Cell(5,4)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "34÷5=6, 4"

$cell = $t.Cell(5,5)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "65÷6=10, 5"

$cell = $t.Cell(9,1)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "52÷8=6, 4"

$cell = $t.Cell(9,2)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "58÷6=9, 4"

$cell = $t.Cell(9,3)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "32÷6=5, 2"

$cell = $t.Cell(9,4)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "27÷2=13, 1"

$cell = $t.Cell(9,5)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "86÷2=43, 0"

$cell = $t.Cell(13,1)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "57÷5=11, 2"

$cell = $t.Cell(13,2)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "67÷2=33, 1"

$cell = $t.Cell(13,3)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "65÷9=7, 2"

$cell = $t.Cell(13,4)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "41÷8=5, 1"

$cell = $t.Cell(13,5)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "42÷5=8, 2"

$cell = $t.Cell(17,1)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "33÷8=4, 1"

$cell = $t.Cell(17,2)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "40÷8=5, 0"

$cell = $t.Cell(17,3)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "77÷8=9, 5"

$cell = $t.Cell(17,4)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "92÷3=30, 2"

$cell = $t.Cell(17,5)
$rng = $cell.Range
$rng2 = $d.Range($rng.Start, $rng.End - 1)
$rng2.Text = "79÷7=11, 2"
